$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F, matching the style used by the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row timestamps for the new time_taken column (values stored as text, like the rest of the data columns)
$timestamps = @(
    "2021-10-05 13:39:47.888700",
    "2021-10-05 13:39:47.888717",
    "2021-10-05 13:39:47.888723",
    "2021-10-05 13:39:47.888729",
    "2021-10-05 13:39:47.888736",
    "2021-10-05 13:39:47.888741",
    "2021-10-05 13:39:47.888747",
    "2021-10-05 13:39:47.888752",
    "2021-10-05 13:39:47.888759",
    "2021-10-05 13:39:47.888765",
    "2021-10-05 13:39:47.888771",
    "2021-10-05 13:39:47.888776",
    "2021-10-05 13:39:47.888781",
    "2021-10-05 13:39:47.888786"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
